$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking scrape refresh: update Price (D) and Volume(1h) (E) for most rows,
# and fix two rows (23/24 and 46/47) where the ranking order changed so the
# coin name/link/price/volume moved to a different row.
#
# Column D holds plain text such as "70.131.83" or "34.00" (note the trailing
# zero, and the use of "." as a thousands separator) which Excel would silently
# reinterpret as a number if assigned directly, losing the original formatting.
# Force Text format, assign the literal string, then restore the default "Normal"
# cell style so no stray formatting is left behind on cells that were unstyled.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '70.131.83'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.39%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.615.36'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +3.33%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '605.35'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.41%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '195.61'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.86%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.628'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.50%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  -1.55%  '

$ws.Range("E10").Value = '  -0.06%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '54.07'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.14%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0000306'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.24%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '9.57'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.06%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.191.58'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +3.29%  '

$ws.Range("E15").Value = '  +4.99%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '592.82'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.66%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '19.22'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.80%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '70.290.58'
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '3.618.52'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +3.51%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.122'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.63%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.995'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.72%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '17.83'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.97%  '

$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.17'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +3.11%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '103.06'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.33%  '

$ws.Range("E25").Value = '  +1.07%  '

$ws.Range("E26").Value = '  -1.17%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.78'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.69%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '9.62'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.14%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '34.00'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.23%  '

$ws.Range("E30").Value = '  -1.09%  '

$ws.Range("E31").Value = '  -0.93%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '12.36'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.88%  '

$ws.Range("E33").Value = '  +1.41%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '63.24'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.60%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0₃0894'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +10.92%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '3.943.60'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +5.71%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.19'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +6.81%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '526.53'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +3.19%  '

$ws.Range("E39").Value = '  +0.13%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '37.37'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.10%  '

$ws.Range("E41").Value = '  +0.95%  '

$ws.Range("E42").Value = '  +1.16%  '

$ws.Range("E43").Value = '  -1.86%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0457'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.20%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.87'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.56%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.37'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.48%  '

$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.141'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.08%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '8.64'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.17%  '

$ws.Range("E49").Value = '  -0.14%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.000253'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +5.34%  '

$ws.Range("E51").Value = '  +3.52%  '
